$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 15: Database Design meeting mapping entities to user views --
# (written first so the new shared strings land in the same order the
#  original author's save produced)
$ws.Range("B15").Value = "12.11."
$ws.Range("D15").Value = "Mapping Entities <-> User Views"
$ws.Range("C15").Value = "Database Design"

# --- New dates for the "user description" / "user views" rows -----------
$ws.Range("B8").Value = "05.10."
$ws.Range("B9").Value = "05.10."

# --- "user view review" row now also carries a date (18.10.) in B -------
$ws.Range("B10").Value = "18.10."

# --- Header row: drop the "Datum" column header (column E) ---------------
$ws.Range("E6").ClearContents()

# --- Column E no longer holds meeting dates, clear the rest too ---------
$ws.Range("E7:E10").ClearContents()

# --- Rows 11-14 are unchanged in content/position ------------------------

# --- Update selection to match the saved workbook -------------------------
$ws.Range("E6:E10").Select()
